# Add a "skype" column to the add_institute_details sheet, populate it with
# skype-id hyperlinks, and make that sheet the active tab (per commit:
# "added skype id in institute management and updated test data and xpaths").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("add_institute_details")

# Insert a new column at F — the previous "runmode" column (F) shifts to G,
# leaving a blank F column ready for the new "skype" data.
[void]$ws.Columns("F").Insert()

# Keep the new column's width in line with its neighbour (E).
$neighborWidth = $ws.Columns("E").ColumnWidth()
$ws.Columns("F").ColumnWidth = $neighborWidth
$ws.Columns("E").ColumnWidth = $neighborWidth

# Header cell.
$ws.Range("F1").NumberFormat = "@"
$ws.Range("F1").Value = "skype"

# Data rows — store each institute's skype id as a hyperlink (mirroring how
# the existing "email" column (C) stores mailto hyperlinks), formatted as text.
$ws.Range("F2").NumberFormat = "@"
[void]$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:sliit@skype", "", "", "sliit@skype")

$ws.Range("F3").NumberFormat = "@"
[void]$ws.Hyperlinks.Add($ws.Range("F3"), "mailto:nibm@skype", "", "", "nibm@skype")

$ws.Range("F4").NumberFormat = "@"
[void]$ws.Hyperlinks.Add($ws.Range("F4"), "mailto:idm@skype", "", "", "idm@skype")

$ws.Range("F5").NumberFormat = "@"
[void]$ws.Hyperlinks.Add($ws.Range("F5"), "mailto:acbt@skype", "", "", "acbt@skype")

# Make add_institute_details the active sheet/selection, matching the
# updated workbook + sheet view state.
[void]$ws.Activate()
[void]$ws.Range("G5").Select()
